# Apply the "Working Deep Learning Algorithm" update:
#  - Sheet2 ("Sequence"): append rows 11-17 with new topic labels
#  - Sheet1: add a new hyperlinked row (B5) pointing to the weight
#    initialisation article referenced from the new "Xavier" topic

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sequence")

# --- Sheet2: new sequence entries (first batch) ---
$ws2.Range("A11").Value = "Activation Functon"
$ws2.Range("A12").Value = "Variants of Relu"
$ws2.Range("A13").Value = "Weight Initialisation Technique"

# --- Sheet1: new hyperlink row (authored in between the two Sheet2 batches) ---
$ws1.Hyperlinks.Add($ws1.Range("B5"), "https://www.deeplearning.ai/ai-notes/initialization/index.html", [System.Type]::Missing, [System.Type]::Missing, "https://www.deeplearning.ai/ai-notes/initialization/index.html")
$ws1.Range("B5").Style = "Hyperlink"

# --- Sheet2: new sequence entries (second batch) ---
$ws2.Range("A14").Value = "Xavier"
$ws2.Range("A15").Value = "Batch Normalisation"
$ws2.Range("A16").Value = "Optimiser"
$ws2.Range("A17").Value = "Exponentially Weighted Moving Average"

# --- Update selections to match where the author left off editing ---
$ws1.Range("B8").Select() | Out-Null
$ws2.Range("A19").Select() | Out-Null
